$wb = $excel.ActiveWorkbook

# --- idsw.datafetch.core: add SharePointDownloader class rows ---
$ws1 = $wb.Worksheets.Item("idsw.datafetch.core")

$ws1.Cells.Item(54, 1).Value = 53
$ws1.Cells.Item(54, 2).Value = "SharePointDownloader"
$ws1.Cells.Item(54, 3).Value = "get_token"

$ws1.Cells.Item(55, 1).Value = 54
$ws1.Cells.Item(55, 2).Value = "SharePointDownloader"
$ws1.Cells.Item(55, 3).Value = "get_response_id"

$ws1.Cells.Item(56, 1).Value = 55
$ws1.Cells.Item(56, 2).Value = "SharePointDownloader"
$ws1.Cells.Item(56, 3).Value = "get_drive_id"

$ws1.Cells.Item(57, 1).Value = 56
$ws1.Cells.Item(57, 2).Value = "SharePointDownloader"
$ws1.Cells.Item(57, 3).Value = "find_file"

$ws1.Cells.Item(58, 1).Value = 57
$ws1.Cells.Item(58, 2).Value = "SharePointDownloader"
$ws1.Cells.Item(58, 3).Value = "download_file"

# --- idsw.modelling.nonsupervised: add benford_outliers_detection row ---
$ws2 = $wb.Worksheets.Item("idsw.modelling.nonsupervised")

$ws2.Cells.Item(6, 1).Value = 5
$ws2.Cells.Item(6, 3).Value = "benford_outliers_detection"
